$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# Drop the "Readme" sheet and the old "Change rates" sheet (superseded by
# "Change rates_v2"); keep "Reference countries" and rename
# "Change rates_v2" to "Change rates" so only the needed sheets remain.
$wb.Worksheets("Readme").Delete()
$wb.Worksheets("Change rates").Delete()
$wb.Worksheets("Change rates_v2").Name = "Change rates"

# Make "Reference countries" the active/selected sheet again, with the
# cursor parked just below the used range.
$ws = $wb.Worksheets("Reference countries")
$ws.Activate()
$ws.Range("B10").Select()
